$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the existing age-band data (rows 2-4), which
# pushes the former row 2 ("15-19 years" ...) down to row 5, etc.
$ws.Range("A2:A4").EntireRow.Insert()

# New "stage distribution" rows with the three new age-bucket labels and
# zeroed out values for all three series (Male and female / Male / Female).
$ws.Range("A2").Value = "(0,4]"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("A3").Value = "(4,9]"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("A4").Value = "(9,14]"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

# Update the chart's series formulas so they point at the data's new
# location (rows 5-19 instead of 2-16).
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$s1 = $chart.FullSeriesCollection(1)
$s1.Formula = "=SERIES(""Male and female"",'Sheet 1'!`$A`$5:`$A`$19,'Sheet 1'!`$B`$5:`$B`$19,1)"

$s2 = $chart.FullSeriesCollection(2)
$s2.Formula = "=SERIES(""Male"",'Sheet 1'!`$A`$5:`$A`$19,'Sheet 1'!`$C`$5:`$C`$19,2)"

$s3 = $chart.FullSeriesCollection(3)
$s3.Formula = "=SERIES(""Female"",'Sheet 1'!`$A`$5:`$A`$19,'Sheet 1'!`$D`$5:`$D`$19,3)"

# The chart graphic itself keeps its size/left position, but slides down by
# the height of the three newly inserted rows so it stays anchored below
# the data table.
$co.Top = $co.Top + 45

# Match the author's last active selection.
$ws.Range("F5").Select()
